# Update the crawl timestamp (column O) for every data row from the old
# crawl time to the new one, and flag two products that went out of stock
# online by rewriting their aria-label text (column M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-12-28 12:55:15"
$newTimestamp = "2022-12-28 20:49:41"

$lastRow = 398
$timestampCol = 15   # column O

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timestampCol)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Row 192: "Prix Garantie Croissants 6 Stück" is now marked as unavailable online.
$ws.Range("M192").Value = "Prix Garantie Croissants 6 Stück - Online kein Bestand 1.95 Schweizer Franken"

# Row 261: "Pasquier Schokobrötchen 8St" is now marked as unavailable online.
$ws.Range("M261").Value = "Pasquier Schokobrötchen 8St - Online kein Bestand 4.60 Schweizer Franken"
